$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new text value, taken from the
# refreshed "cryptos" price/volume snapshot (GitHub Actions data refresh).
$updates = @(
    @{Cell = "D2"; Value = "29.348.91"},
    @{Cell = "E2"; Value = "  -0.26%  "},
    @{Cell = "D3"; Value = "1.843.90"},
    @{Cell = "E3"; Value = "  -0.35%  "},
    @{Cell = "D4"; Value = "0.9999"},
    @{Cell = "E4"; Value = "  -0.03%  "},
    @{Cell = "D5"; Value = "240.83"},
    @{Cell = "E5"; Value = "  +0.19%  "},
    @{Cell = "D6"; Value = "0.6288"},
    @{Cell = "E6"; Value = "  -0.14%  "},
    @{Cell = "D7"; Value = "1.001"},
    @{Cell = "E7"; Value = "  -0.05%  "},
    @{Cell = "D8"; Value = "0.07471"},
    @{Cell = "E8"; Value = "  -2.34%  "},
    @{Cell = "D9"; Value = "0.2900"},
    @{Cell = "E9"; Value = "  -0.49%  "},
    @{Cell = "D10"; Value = "24.35"},
    @{Cell = "E10"; Value = "  -2.26%  "},
    @{Cell = "D11"; Value = "0.07736"},
    @{Cell = "E11"; Value = "  -0.13%  "},
    @{Cell = "D12"; Value = "1.843.76"},
    @{Cell = "E12"; Value = "  -2.35%  "},
    @{Cell = "D13"; Value = "5.000"},
    @{Cell = "E13"; Value = "  -0.69%  "},
    @{Cell = "D14"; Value = "0.6780"},
    @{Cell = "E14"; Value = "  -0.48%  "},
    @{Cell = "D15"; Value = "0.00001022"},
    @{Cell = "E15"; Value = "  -4.09%  "},
    @{Cell = "D16"; Value = "82.02"},
    @{Cell = "E16"; Value = "  -1.63%  "},
    @{Cell = "D17"; Value = "6.127"},
    @{Cell = "E17"; Value = "  -0.89%  "},
    @{Cell = "D18"; Value = "29.390.02"},
    @{Cell = "E18"; Value = "  -0.38%  "},
    @{Cell = "D19"; Value = "227.93"},
    @{Cell = "E19"; Value = "  -0.28%  "},
    @{Cell = "D20"; Value = "12.30"},
    @{Cell = "E20"; Value = "  -0.31%  "},
    @{Cell = "D21"; Value = "1.000"},
    @{Cell = "E21"; Value = "  -0.08%  "},
    @{Cell = "D22"; Value = "7.430"},
    @{Cell = "E22"; Value = "  -0.51%  "},
    @{Cell = "D23"; Value = "1.001"},
    @{Cell = "E23"; Value = "  +0.04%  "},
    @{Cell = "D24"; Value = "158.60"},
    @{Cell = "E24"; Value = "  +0.41%  "},
    @{Cell = "D25"; Value = "0.1374"},
    @{Cell = "E25"; Value = "  -0.85%  "},
    @{Cell = "D26"; Value = "8.422"},
    @{Cell = "E26"; Value = "  -0.30%  "},
    @{Cell = "D27"; Value = "17.54"},
    @{Cell = "E27"; Value = "  -0.91%  "},
    @{Cell = "D28"; Value = "0.06390"},
    @{Cell = "E28"; Value = "  +13.62%  "},
    @{Cell = "D29"; Value = "1.393"},
    @{Cell = "E29"; Value = "  +1.09%  "},
    @{Cell = "E30"; Value = "  +0.51%  "},
    @{Cell = "D31"; Value = "4.091"},
    @{Cell = "E31"; Value = "  -1.02%  "},
    @{Cell = "D32"; Value = "4.048"},
    @{Cell = "E33"; Value = "  -1.45%  "},
    @{Cell = "D35"; Value = "0.6962"},
    @{Cell = "E35"; Value = "  -0.91%  "},
    @{Cell = "D36"; Value = "2.585"},
    @{Cell = "E36"; Value = "  -0.29%  "},
    @{Cell = "D37"; Value = "1.265.88"},
    @{Cell = "E37"; Value = "  +2.97%  "},
    @{Cell = "E38"; Value = "  +4.20%  "},
    @{Cell = "D39"; Value = "0.01814"},
    @{Cell = "E39"; Value = "  +0.48%  "},
    @{Cell = "D40"; Value = "6.542"},
    @{Cell = "E40"; Value = "  +1.49%  "},
    @{Cell = "D41"; Value = "0.9095"},
    @{Cell = "E41"; Value = "  +0.01%  "},
    @{Cell = "D42"; Value = "0.9996"},
    @{Cell = "E42"; Value = "  -0.16%  "},
    @{Cell = "D43"; Value = "2.004.63"},
    @{Cell = "E43"; Value = "  -12.00%  "},
    @{Cell = "D44"; Value = "101.31"},
    @{Cell = "E44"; Value = "  -1.11%  "},
    @{Cell = "D45"; Value = "66.31"},
    @{Cell = "E45"; Value = "  +0.32%  "},
    @{Cell = "D46"; Value = "0.1172"},
    @{Cell = "E46"; Value = "  +1.43%  "},
    @{Cell = "D47"; Value = "7.030"},
    @{Cell = "E47"; Value = "  -2.31%  "},
    @{Cell = "B48"; Value = "EnergySwap"},
    @{Cell = "C48"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"},
    @{Cell = "D48"; Value = "9.045"},
    @{Cell = "E48"; Value = "  +0.69%  "},
    @{Cell = "B49"; Value = "BabyDogeCoin"},
    @{Cell = "C49"; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"},
    @{Cell = "D49"; Value = "0.00000000115"},
    @{Cell = "E49"; Value = "  -4.16%  "},
    @{Cell = "D51"; Value = "1.673"},
    @{Cell = "E51"; Value = "  -0.73%  "}
)

foreach ($item in $updates) {
    $cell = $ws.Range($item.Cell)
    # Prefix with an apostrophe so Excel keeps numeric-looking strings
    # (e.g. "1.001", "5.000") as plain text instead of coercing them to
    # numbers, matching the source data which is all text.
    $cell.Value = "'" + $item.Value
    # Clear the resulting quote-prefix formatting so the cell style stays
    # the same as before the edit.
    $cell.Style = "Normal"
}

Write-Output ("Updated {0} cells" -f $updates.Count)
